# Apply cell value updates to Sheet1 as described in the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 17).Value = 1.71
$ws.Cells.Item(2, 19).Value = 2.8
$ws.Cells.Item(2, 32).Value = 19.5
$ws.Cells.Item(3, 6).Value = 1.52
$ws.Cells.Item(3, 16).Value = 1.83
$ws.Cells.Item(3, 17).Value = 1.86
$ws.Cells.Item(3, 20).Value = 2.02
$ws.Cells.Item(3, 21).Value = 1.82
$ws.Cells.Item(4, 20).Value = 1.94
$ws.Cells.Item(4, 21).Value = 1.86
$ws.Cells.Item(5, 6).Value = 3.9
$ws.Cells.Item(5, 7).Value = 5.2
$ws.Cells.Item(5, 8).Value = 1.9
$ws.Cells.Item(5, 11).Value = 4.1
$ws.Cells.Item(5, 17).Value = 1.73
$ws.Cells.Item(5, 18).Value = 1.37
$ws.Cells.Item(5, 19).Value = 3.15
$ws.Cells.Item(5, 20).Value = 1.73
$ws.Cells.Item(5, 21).Value = 2.08
$ws.Cells.Item(5, 23).Value = 1.27
$ws.Cells.Item(6, 8).Value = 1.41
$ws.Cells.Item(6, 9).Value = 1.44
$ws.Cells.Item(6, 11).Value = 5.5
$ws.Cells.Item(6, 16).Value = 2.26
$ws.Cells.Item(9, 19).Value = 3.5
$ws.Cells.Item(10, 7).Value = 4.4
$ws.Cells.Item(10, 8).Value = 2.22
$ws.Cells.Item(10, 9).Value = 2.7
$ws.Cells.Item(10, 10).Value = 2.9
$ws.Cells.Item(10, 11).Value = 3.75
$ws.Cells.Item(10, 14).Value = 1.81
$ws.Cells.Item(10, 16).Value = 1.81
$ws.Cells.Item(10, 17).Value = 1.88
$ws.Cells.Item(10, 22).Value = 1.59
$ws.Cells.Item(10, 23).Value = 1.29
$ws.Cells.Item(10, 25).Value = 970
$ws.Cells.Item(10, 26).Value = 980
$ws.Cells.Item(10, 27).Value = 980
$ws.Cells.Item(10, 29).Value = 11
$ws.Cells.Item(10, 31).Value = 980
$ws.Cells.Item(10, 32).Value = 980
$ws.Cells.Item(10, 33).Value = 980
$ws.Cells.Item(10, 34).Value = 980
$ws.Cells.Item(10, 35).Value = 65
$ws.Cells.Item(10, 37).Value = 65
$ws.Cells.Item(11, 6).Value = 1.92
$ws.Cells.Item(11, 7).Value = 2.08
$ws.Cells.Item(11, 8).Value = 3.8
$ws.Cells.Item(11, 9).Value = 4.9
$ws.Cells.Item(11, 11).Value = 4.2
$ws.Cells.Item(11, 18).Value = 1.37
$ws.Cells.Item(11, 19).Value = 2.56
$ws.Cells.Item(11, 22).Value = 1.29
$ws.Cells.Item(11, 23).Value = 1.92
$ws.Cells.Item(11, 24).Value = 26
$ws.Cells.Item(11, 25).Value = 25
$ws.Cells.Item(11, 26).Value = 46
$ws.Cells.Item(11, 27).Value = 100
$ws.Cells.Item(11, 28).Value = 15.5
$ws.Cells.Item(11, 29).Value = 13
$ws.Cells.Item(11, 30).Value = 24
$ws.Cells.Item(11, 31).Value = 70
$ws.Cells.Item(11, 32).Value = 20
$ws.Cells.Item(11, 33).Value = 15.5
$ws.Cells.Item(11, 34).Value = 25
$ws.Cells.Item(11, 35).Value = 70
$ws.Cells.Item(11, 36).Value = 34
$ws.Cells.Item(11, 37).Value = 29
$ws.Cells.Item(11, 38).Value = 48
$ws.Cells.Item(12, 17).Value = 1.4
$ws.Cells.Item(12, 18).Value = 1.88
$ws.Cells.Item(12, 19).Value = 1.94
$ws.Cells.Item(13, 12).Value = 1.24
$ws.Cells.Item(13, 18).Value = 1.73
$ws.Cells.Item(13, 19).Value = 2
$ws.Cells.Item(14, 6).Value = 1.86
$ws.Cells.Item(14, 9).Value = 4.6
$ws.Cells.Item(14, 11).Value = 5
$ws.Cells.Item(14, 13).Value = 1.03
$ws.Cells.Item(14, 14).Value = 5.3
$ws.Cells.Item(14, 18).Value = 1.59
$ws.Cells.Item(14, 19).Value = 2.42
$ws.Cells.Item(14, 20).Value = 1.58
$ws.Cells.Item(14, 21).Value = 2.46
$ws.Cells.Item(14, 22).Value = 1.3
$ws.Cells.Item(14, 24).Value = 30
$ws.Cells.Item(14, 25).Value = 25
$ws.Cells.Item(14, 26).Value = 42
$ws.Cells.Item(14, 27).Value = 95
$ws.Cells.Item(14, 28).Value = 16
$ws.Cells.Item(14, 29).Value = 12.5
$ws.Cells.Item(14, 30).Value = 21
$ws.Cells.Item(14, 31).Value = 55
$ws.Cells.Item(14, 32).Value = 17.5
$ws.Cells.Item(14, 33).Value = 13
$ws.Cells.Item(14, 34).Value = 19.5
$ws.Cells.Item(14, 35).Value = 55
$ws.Cells.Item(14, 36).Value = 27
$ws.Cells.Item(14, 37).Value = 22
$ws.Cells.Item(14, 38).Value = 34
$ws.Cells.Item(14, 39).Value = 80
$ws.Cells.Item(14, 40).Value = 11
$ws.Cells.Item(14, 41).Value = 38
$ws.Cells.Item(15, 6).Value = 1.04
$ws.Cells.Item(15, 7).Value = 1.84
$ws.Cells.Item(15, 10).Value = 3.95
$ws.Cells.Item(16, 12).Value = 1.22
$ws.Cells.Item(16, 20).Value = 1.47
$ws.Cells.Item(17, 18).Value = 1.77
$ws.Cells.Item(17, 19).Value = 2.04
$ws.Cells.Item(18, 12).Value = 1.34
$ws.Cells.Item(18, 18).Value = 1.16
$ws.Cells.Item(18, 24).Value = 15.5
$ws.Cells.Item(18, 25).Value = 17.5
$ws.Cells.Item(18, 26).Value = 40
$ws.Cells.Item(18, 28).Value = 11.5
$ws.Cells.Item(18, 29).Value = 10.5
$ws.Cells.Item(18, 30).Value = 25
$ws.Cells.Item(18, 31).Value = 90
$ws.Cells.Item(18, 32).Value = 19
$ws.Cells.Item(18, 33).Value = 16
$ws.Cells.Item(18, 34).Value = 30
$ws.Cells.Item(18, 35).Value = 100
$ws.Cells.Item(18, 36).Value = 44
$ws.Cells.Item(18, 37).Value = 40
$ws.Cells.Item(18, 38).Value = 70
$ws.Cells.Item(20, 12).Value = 1.48
$ws.Cells.Item(20, 14).Value = 1.01
$ws.Cells.Item(20, 15).Value = 1.52
$ws.Cells.Item(20, 16).Value = 1.08
$ws.Cells.Item(20, 17).Value = 1.52
$ws.Cells.Item(20, 19).Value = 1.01
$ws.Cells.Item(24, 6).Value = 2.6
$ws.Cells.Item(24, 7).Value = 2.8
$ws.Cells.Item(24, 8).Value = 3.55
$ws.Cells.Item(24, 9).Value = 3.9
$ws.Cells.Item(24, 10).Value = 2.68
$ws.Cells.Item(24, 11).Value = 2.9
$ws.Cells.Item(24, 13).Value = 1.19
$ws.Cells.Item(24, 14).Value = 2.08
$ws.Cells.Item(24, 15).Value = 1.81
$ws.Cells.Item(24, 16).Value = 1.34
$ws.Cells.Item(24, 17).Value = 3.5
$ws.Cells.Item(24, 18).Value = 1.11
$ws.Cells.Item(24, 19).Value = 7.2
$ws.Cells.Item(24, 20).Value = 2.52
$ws.Cells.Item(24, 21).Value = 1.56
$ws.Cells.Item(24, 22).Value = 1.34
$ws.Cells.Item(24, 23).Value = 1.55
$ws.Cells.Item(24, 33).Value = 970
$ws.Cells.Item(28, 6).Value = 2.32
$ws.Cells.Item(28, 11).Value = 5.5
